$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 258; this shifts existing rows 258-268
# down to 259-269 and extends the used range to R269.
$ws.Rows.Item(258).Insert()

# Populate the newly inserted row 258 with the new weekly record.
$ws.Range("A258").Value = 10
$ws.Range("B258").Value = "Vega Modelo de Temuco"
$ws.Range("C258").Value = "La Araucanía"
$ws.Range("D258").Value = 45041
$ws.Range("E258").Value = 9
$ws.Range("F258").Value = 100112012
$ws.Range("G258").Value = "Espinaca"
$ws.Range("H258").Value = "Sin especificar"
$ws.Range("I258").Value = "Primera"
$ws.Range("J258").Value = 40
$ws.Range("K258").Value = 10000
$ws.Range("L258").Value = 10000
$ws.Range("M258").Value = 10000
$ws.Range("N258").Value = "$/docena de atados"
$ws.Range("O258").Value = "Región de La Araucanía"
$ws.Range("P258").Value = 3333
$ws.Range("Q258").Value = 3
$ws.Range("R258").Value = "Hortaliza"
